$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New chapter counts per year (row 2..27 -> years 1997..2022)
$values = @(8, 11, 22, 24, 25, 26, 41, 31, 31, 35, 35, 34, 31, 33, 49, 48, 50, 55, 53, 52, 52, 49, 61, 48, 59, 19)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 2).Value = $values[$i]
}

# B2 previously used a distinct date-like number format (style index 5); align it
# with the rest of the column (style used by B3) now that it's a plain number.
$ws.Range("B3").Copy()
$ws.Range("B2").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Update the view: zoom and selection
$excel.ActiveWindow.Zoom = 99
$ws.Range("D13").Select()
